# Generate Report for Handoff
#
# A new handoff xliff generation run updated the "Latest Handoff
# Datetime" / "Latest HO Xliff Generate Date" timestamps for the file
# 0e4d11ec-e77e-4851-9829-212f3bba15fe.md (row 5 on every sheet).

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-08-17 10:39:18"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-08-17 10:39:13"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-08-17 10:39:18"
